$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2807076.42
$ws.Range("C7").Value = -36.821357982335
$ws.Range("D7").Value = 2866
$ws.Range("E7").Value = 2866
$ws.Range("F7").Value = 979.4404815073273
$ws.Range("G7").Value = 4.401273061989319
